# Actualización desde MV -datos-
# Updates the daily currency-parity table: fixes a bad value in D190 and
# appends two new date rows (04-10-2021 and 05-10-2021) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the erroneous value in D190 (was a mis-scaled outlier, should be the
#    same order of magnitude as the surrounding "Bolívar Venezolano" column).
# ---------------------------------------------------------------------------
$ws.Cells.Item(190, 4).Value = 4.1283

# ---------------------------------------------------------------------------
# Helper: write a date-looking label into column A as literal text (not an
# auto-converted serial date). We build it as a formula returning the text,
# then Copy/PasteSpecial(values-only) to "flatten" it back into a plain
# shared-string cell, leaving cell formatting untouched.
# ---------------------------------------------------------------------------
function Set-TextLabel {
    param($row, $text)
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

# ---------------------------------------------------------------------------
# Helper: write an entire B:BH row of numeric data in one shot via a 2D
# SAFEARRAY, matching how Excel's Range.Value setter expects multi-cell data.
# ---------------------------------------------------------------------------
function Set-RowValues {
    param($row, $values)
    $arr = New-Object 'object[,]' 1, $values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $rng = $ws.Range(("B{0}:BH{0}" -f $row))
    $rng.Value = $arr
}

# ---------------------------------------------------------------------------
# 2) Row 191 — 04-10-2021
# ---------------------------------------------------------------------------
Set-TextLabel 191 "04-10-2021"
Set-RowValues 191 @(
    33.625, 1, 4.1283, 6.86, 624.45, 21.8209, 6.4116, 129.46, 8.6213, 8.7472,
    0.7084, 3.672, 9.0548, 1.3748, 1.264, 1, 2.1124, 0.82, 1, 1.3549, 7.7849,
    1.4401, 27.837, 22737, 0.8622, 307.84, 102.6012, 0.9297, 6902, 26.61,
    4.2618, 15.69, 0.7374, 8.8528, 4.1277, 98.7721, 803.9, 3783.96, 24,
    56.231, 50.79, 20.4288, 42.78, 7.7283, 14.858, 5.3652, 42000, 3.75,
    4.177, 72.6958, 74.12, 14305, 169.75, 3.215, 426.14, 1187.73, 110.99,
    6.4302, 3.9469
)

# ---------------------------------------------------------------------------
# 3) Row 192 — 05-10-2021
# ---------------------------------------------------------------------------
Set-TextLabel 192 "05-10-2021"
Set-RowValues 192 @(
    33.749, 1, 4.1713, 6.86, 623.53, 21.8194, 6.4008, 127.98, 8.582, 8.7275,
    0.7076, 3.672, 9.0028, 1.3729, 1.258, 1, 2.107, 0.82, 1, 1.3571, 7.7868,
    1.4362, 27.863, 22734, 0.8607, 306.58, 102.3917, 0.9249, 6905.4, 26.3598,
    4.2583, 15.66, 0.7345, 8.8609, 4.13, 98.8767, 805.89, 3789, 24, 56.217,
    50.68, 20.5513, 42.8339, 7.7287, 15.0611, 5.4357, 42000, 3.75, 4.173,
    72.4712, 74.3075, 14265, 170.5, 3.2236, 424.9, 1181.97, 110.9, 6.4494,
    3.9533
)

$excel.CutCopyMode = $false
